# feat(ExporterEffects): complete export of assessments grid
#
# - Rename the "SURA" sheet to "Assessments"
# - Main sheet: add a summary row for the Assessments grid (label +
#   COUNTA formula), matching the existing "Conversations" summary row,
#   and widen column A to fit the new, longer label
# - Make "Main" the active sheet/tab (it was "Conversations" before)
# - Leave a selection parked on G20 of the renamed "Assessments" sheet

$wb = $excel.ActiveWorkbook

$mainSheet = $wb.Worksheets.Item("Main")
$conversationsSheet = $wb.Worksheets.Item("Conversations")
$assessmentsSheet = $wb.Worksheets.Item("SURA")

# 1. Rename SURA -> Assessments (rId/sheetId stay the same, only the name changes)
$assessmentsSheet.Name = "Assessments"

# 2. Main: new row 3 mirroring row 2's "Conversations" summary, but for Assessments
$mainSheet.Range("A3").Value = "Assessments"
$mainSheet.Range("B3").Formula = "=COUNTA(Assessments!A:A)"

# 3. Main: column A needs to be wide enough for the new "Assessments" label
$mainSheet.Columns.Item(1).ColumnWidth = 20.33

# 4. Leave a selection on the renamed Assessments sheet at G20
[void]$assessmentsSheet.Range("G20").Select()

# 5. Make Main the active sheet/tab, with A4 selected (row below the new data)
[void]$mainSheet.Activate()
[void]$mainSheet.Range("A4").Select()
